# Cronograma.xlsx update: inserting new schedule items into the "Projeto" sheet,
# per commit "Atualização da Planilha de Cronograma, inserindo novos itens
# (faltam ajustes nas macros)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projeto")

# ---------------------------------------------------------------------------
# 1) Make room: stretch the generic activity-row style (B:G of row 9) down
#    across rows 10-26 so every row - old and new - carries the same cell
#    formatting (text style for B, numeric style for C-F, percent style G)
#    before we touch any values.
# ---------------------------------------------------------------------------
$ws.Range("B9:G9").Copy($ws.Range("B10:G26"))

# Row 17 ("Montagem maquete") keeps the slightly different D/F font that the
# old row 13 used to carry (style id 21 instead of 15) - grab it from there
# before row 13 gets overwritten below.
$ws.Range("D13:D13").Copy($ws.Range("D17"))
$ws.Range("F13:F13").Copy($ws.Range("F17"))

# ---------------------------------------------------------------------------
# 2) Row heights: rows 9-16 use 18.75, rows 17-26 use 19.
# ---------------------------------------------------------------------------
$ws.Range("9:16").RowHeight = 18.75
$ws.Range("17:26").RowHeight = 19

# ---------------------------------------------------------------------------
# 3) Write the final activity table, row by row.
# ---------------------------------------------------------------------------
$rows = @(
    @{Row=9;  B='Plano de projeto';                                 C=1;  D=5;  E=1;    F=5;    G=1},
    @{Row=10; B='Blog - Criação e Atualização';                     C=2;  D=35; E=$null; F=$null; G=1},
    @{Row=11; B='GitHub - Criação e Atualização';                   C=3;  D=35; E=$null; F=$null; G=1},
    @{Row=12; B='Casos de uso';                                     C=5;  D=3;  E=5;    F=4;    G=1},
    @{Row=13; B='Outros digramas UML';                              C=6;  D=8;  E=6;    F=$null; G=$null},
    @{Row=14; B='Preparação da Entrega 1';                          C=7;  D=1;  E=7;    F=1;    G=$null},
    @{Row=15; B='Preparação da Entrega 2';                          C=14; D=1;  E=14;   F=1;    G=$null},
    @{Row=16; B='Preparação da Entrega 3';                          C=21; D=1;  E=21;   F=1;    G=$null},
    @{Row=17; B='Montagem maquete';                                 C=15; D=2;  E=$null; F=$null; G=0.35},
    @{Row=18; B='Testes dos equipamentos';                          C=8;  D=4;  E=8;    F=2;    G=0.2},
    @{Row=19; B='Preparação e Montagem inicial dos dispositivos';   C=9;  D=2;  E=9;    F=2;    G=$null},
    @{Row=20; B='Inicio da programação / Blynk / Circuitos';        C=12; D=5;  E=$null; F=$null; G=0},
    @{Row=21; B='Programação / Blynk / Adaptação dos Circuitos';    C=17; D=11; E=$null; F=$null; G=0},
    @{Row=22; B='Testes';                                           C=27; D=2;  E=$null; F=$null; G=$null},
    @{Row=23; B='Ajustes e correções';                              C=29; D=4;  E=$null; F=$null; G=0},
    @{Row=24; B='Testes finais';                                    C=33; D=2;  E=$null; F=$null; G=0},
    @{Row=25; B='Preparação da apresentação ';                      C=29; D=6;  E=$null; F=$null; G=0},
    @{Row=26; B='Documentação final';                                C=32; D=6;  E=$null; F=$null; G=0}
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 2).Value = $r.B
    if ($null -eq $r.C) { $ws.Cells.Item($n, 3).Value = "" } else { $ws.Cells.Item($n, 3).Value = $r.C }
    if ($null -eq $r.D) { $ws.Cells.Item($n, 4).Value = "" } else { $ws.Cells.Item($n, 4).Value = $r.D }
    if ($null -eq $r.E) { $ws.Cells.Item($n, 5).Value = "" } else { $ws.Cells.Item($n, 5).Value = $r.E }
    if ($null -eq $r.F) { $ws.Cells.Item($n, 6).Value = "" } else { $ws.Cells.Item($n, 6).Value = $r.F }
    if ($null -eq $r.G) { $ws.Cells.Item($n, 7).Value = "" } else { $ws.Cells.Item($n, 7).Value = $r.G }
}

# ---------------------------------------------------------------------------
# 4) "DIAS DO PROJETO" period selector moved from day 10 to day 11, and the
#    active selection follows it (N4 -> N3).
# ---------------------------------------------------------------------------
$ws.Range("N3").Value = 11
$ws.Range("N3").Select()

# ---------------------------------------------------------------------------
# 5) Conditional-formatting ranges grow along with the table: the big
#    activity block now reaches row 26, and the always-on stripe row moves
#    from 21 to 27.
# ---------------------------------------------------------------------------
$bigRangeOld = $ws.Range("I12:AS20")
$bigRangeNew = $ws.Range("I12:AS26")
$bigFcs = $bigRangeOld.FormatConditions
for ($i = 1; $i -le $bigFcs.Count; $i++) {
    $bigFcs.Item($i).ModifyAppliesToRange($bigRangeNew)
}

$stripeFcs = $ws.Range("B21:AS21").FormatConditions
$stripeFcs.Item(1).ModifyAppliesToRange($ws.Range("B27:AS27"))
